$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Fix the product name text (missing hyphen after "342") on both sheets
$ws1.Range("B1").Value = "342-MS-EPP-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"
$ws2.Range("B1").Value = "342-MS-EPP-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"

# Update selection on the input sheet (was A6:B6, now B1) without leaving it active
$ws1.Activate()
$ws1.Range("B1").Select()

# Make the output sheet the active/selected tab with B1 selected
$ws2.Activate()
$ws2.Range("B1").Select()
